# Update the "K" column (column G) values in the save-data sheet.
# This reflects regenerating save_data to use K (strikeouts) instead of Strike#,
# along with the recalculated std/mean and s_vals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 2
    3  = 0
    4  = 0
    5  = 0
    6  = 1
    7  = 0
    8  = 1
    9  = 0
    10 = 1
    11 = 4
    12 = 2
    13 = 1
    14 = 0
    15 = 1
    16 = 2
    17 = 1
    18 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
